$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<other>"
$ws.Range("C2").Value = 50

# Row 3
$ws.Range("B3").Value = "<up>"
$ws.Range("C3").Value = 53

# Row 4
$ws.Range("B4").Value = "<by>"
$ws.Range("C4").Value = 47

# Row 5
$ws.Range("B5").Value = "<it>"
$ws.Range("C5").Value = 48

# Row 6
$ws.Range("B6").Value = "<we>"
$ws.Range("C6").Value = 49

# Row 7
$ws.Range("C7").Value = 47

# Row 8
$ws.Range("B8").Value = "<they>"
$ws.Range("C8").Value = 46

# Row 9
$ws.Range("B9").Value = "<delete>"
$ws.Range("C9").Value = 42

# Row 10
$ws.Range("C10").Value = 47

# Row 11
$ws.Range("C11").Value = 49

# Row 12
$ws.Range("B12").Value = "<them>"
$ws.Range("C12").Value = 48

# Row 13
$ws.Range("B13").Value = "<for>"
$ws.Range("C13").Value = 43

# Row 14
$ws.Range("C14").Value = 50

# Row 15
$ws.Range("B15").Value = "<i>"
$ws.Range("C15").Value = 52

# Row 16
$ws.Range("C16").Value = 49

# Row 17
$ws.Range("C17").Value = 52

# Row 18
$ws.Range("C18").Value = 44
